$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates that are safe as general text (contain multiple
# separators, so Excel will not reinterpret them as numbers)
$ws.Range("D2").Value = "61.795.54"
$ws.Range("D3").Value = "3.395.58"
$ws.Range("D8").Value = "3.394.48"
$ws.Range("D13").Value = "3.973.65"
$ws.Range("D16").Value = "3.401.40"
$ws.Range("D18").Value = "61.815.01"
$ws.Range("D24").Value = "3.532.47"
$ws.Range("D35").Value = "3.426.79"
$ws.Range("D51").Value = "2.363.79"

# Price (column D) updates that look like plain numbers - force the cell to
# Text format first so Excel keeps them as literal strings (matching the
# original inline-string cells) instead of parsing them into numeric values
# (which would silently drop trailing zeros / use scientific notation).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.02"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.57"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.50"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.41"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.20"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.65"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.568"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000128"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.32"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.51"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.42"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.69"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.72"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.66"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.95"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.43"

# Volume(1h) (column E) updates - plain text with padding spaces, never numeric
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("E22").Value = "  +6.56%  "
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("E25").Value = "  +13.22%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  +12.33%  "
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("E51").Value = "  +7.54%  "
